$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 30
$ws.Range("E2").Value = 33.40000152587891
$ws.Range("F2").Value = 35.77999877929688
$ws.Range("G2").Value = 25.17000007629395
$ws.Range("H2").Value = 429709983
$ws.Range("I2").Value = "RXRX"

$ws.Range("D3").Value = 30
$ws.Range("E3").Value = 33.40000152587891
$ws.Range("F3").Value = 35.77999877929688
$ws.Range("G3").Value = 25.17000007629395
$ws.Range("H3").Value = 429709983
$ws.Range("I3").Value = "RXRX"

$ws.Range("D4").Value = 30
$ws.Range("E4").Value = 33.40000152587891
$ws.Range("F4").Value = 35.77999877929688
$ws.Range("G4").Value = 25.17000007629395
$ws.Range("H4").Value = 429709983
$ws.Range("I4").Value = "RXRX"

$ws.Range("D5").Value = 30
$ws.Range("E5").Value = 33.40000152587891
$ws.Range("F5").Value = 35.77999877929688
$ws.Range("G5").Value = 25.17000007629395
$ws.Range("H5").Value = 429709983
$ws.Range("I5").Value = "RXRX"

$ws.Range("D6").Value = 30
$ws.Range("E6").Value = 33.40000152587891
$ws.Range("F6").Value = 35.77999877929688
$ws.Range("G6").Value = 25.17000007629395
$ws.Range("H6").Value = 429709983
$ws.Range("I6").Value = "RXRX"

$ws.Range("D7").Value = 30
$ws.Range("E7").Value = 33.40000152587891
$ws.Range("F7").Value = 35.77999877929688
$ws.Range("G7").Value = 25.17000007629395
$ws.Range("H7").Value = 429709983
$ws.Range("I7").Value = "RXRX"

$ws.Range("D8").Value = 36.72000122070312
$ws.Range("E8").Value = 30.35000038146973
$ws.Range("F8").Value = 42.81000137329102
$ws.Range("G8").Value = 29.20000076293945
$ws.Range("H8").Value = 429709983
$ws.Range("I8").Value = "RXRX"

$ws.Range("D9").Value = 22.96999931335449
$ws.Range("E9").Value = 19.14999961853028
$ws.Range("F9").Value = 23.86000061035156
$ws.Range("G9").Value = 16.5
$ws.Range("H9").Value = 429709983
$ws.Range("I9").Value = "RXRX"

$ws.Range("D10").Value = 17.20000076293945
$ws.Range("E10").Value = 11.84000015258789
$ws.Range("F10").Value = 18.22999954223633
$ws.Range("G10").Value = 9.859999656677246
$ws.Range("H10").Value = 429709983
$ws.Range("I10").Value = "RXRX"

$ws.Range("D11").Value = 7.170000076293945
$ws.Range("E11").Value = 6.199999809265137
$ws.Range("F11").Value = 7.935999870300293
$ws.Range("G11").Value = 5.71999979019165
$ws.Range("H11").Value = 429709983
$ws.Range("I11").Value = "RXRX"

$ws.Range("D12").Value = 8.149999618530273
$ws.Range("E12").Value = 8.460000038146973
$ws.Range("F12").Value = 9.289999961853027
$ws.Range("G12").Value = 7.48199987411499
$ws.Range("H12").Value = 429709983
$ws.Range("I12").Value = "RXRX"

$ws.Range("D13").Value = 10.85000038146973
$ws.Range("E13").Value = 10.55000019073486
$ws.Range("F13").Value = 11.9350004196167
$ws.Range("G13").Value = 9.800000190734863
$ws.Range("H13").Value = 429709983
$ws.Range("I13").Value = "RXRX"

$ws.Range("D14").Value = 7.769999980926514
$ws.Range("E14").Value = 8.329999923706055
$ws.Range("F14").Value = 9.479999542236328
$ws.Range("G14").Value = 7.099999904632568
$ws.Range("H14").Value = 429709983
$ws.Range("I14").Value = "RXRX"

$ws.Range("D15").Value = 6.639999866485596
$ws.Range("E15").Value = 4.769999980926514
$ws.Range("F15").Value = 7
$ws.Range("G15").Value = 4.619999885559082
$ws.Range("H15").Value = 429709983
$ws.Range("I15").Value = "RXRX"

$ws.Range("D16").Value = 7.449999809265137
$ws.Range("E16").Value = 14.11999988555908
$ws.Range("F16").Value = 16.7450008392334
$ws.Range("G16").Value = 6.440000057220459
$ws.Range("H16").Value = 429709983
$ws.Range("I16").Value = "RXRX"

$ws.Range("D17").Value = 7.659999847412109
$ws.Range("E17").Value = 5.28000020980835
$ws.Range("F17").Value = 7.690000057220459
$ws.Range("G17").Value = 5.130000114440918
$ws.Range("H17").Value = 429709983
$ws.Range("I17").Value = "RXRX"

$ws.Range("D18").Value = 9.680000305175779
$ws.Range("E18").Value = 9.409999847412109
$ws.Range("F18").Value = 14.18000030517578
$ws.Range("G18").Value = 9.255000114440918
$ws.Range("H18").Value = 429709983
$ws.Range("I18").Value = "RXRX"

$ws.Range("D19").Value = 9.979999542236328
$ws.Range("E19").Value = 7.820000171661377
$ws.Range("F19").Value = 9.979999542236328
$ws.Range("G19").Value = 7.135000228881836
$ws.Range("H19").Value = 429709983
$ws.Range("I19").Value = "RXRX"

$ws.Range("D20").Value = 7.53000020980835
$ws.Range("E20").Value = 8.199999809265137
$ws.Range("F20").Value = 8.734999656677246
$ws.Range("G20").Value = 7.039999961853027
$ws.Range("H20").Value = 429709983
$ws.Range("I20").Value = "RXRX"

$ws.Range("D21").Value = 6.510000228881836
$ws.Range("E21").Value = 6.320000171661377
$ws.Range("F21").Value = 7.099999904632568
$ws.Range("G21").Value = 5.949999809265137
$ws.Range("H21").Value = 429709983
$ws.Range("I21").Value = "RXRX"

$ws.Range("D22").Value = 6.909999847412109
$ws.Range("E22").Value = 7.239999771118164
$ws.Range("F22").Value = 8.659999847412109
$ws.Range("G22").Value = 6.199999809265137
$ws.Range("H22").Value = 429709983
$ws.Range("I22").Value = "RXRX"

$ws.Range("D23").Value = 5.320000171661377
$ws.Range("E23").Value = 5.590000152587891
$ws.Range("F23").Value = 6.110000133514404
$ws.Range("G23").Value = 3.789999961853027
$ws.Range("H23").Value = 429709983
$ws.Range("I23").Value = "RXRX"

$ws.Range("D24").Value = 5
$ws.Range("E24").Value = 5.949999809265137
$ws.Range("F24").Value = 7.150000095367432
$ws.Range("G24").Value = 4.800000190734863
$ws.Range("H24").Value = 429709983
$ws.Range("I24").Value = "RXRX"
